$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 80.38067699035999
$ws.Range("C2").Value = 18.85621941051061
$ws.Range("D2").Value = 89.25706332453439
$ws.Range("E2").Value = 54.08160640019621
